$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "69.362.93"
$ws.Range("D3").Value = "3.443.84"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "3.435.25"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.193"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.563"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000270"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.997.06"
$ws.Range("D17").Value = "3.453.39"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "581.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "69.391.26"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.846"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "579.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0476"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0954"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "3.236.73"
$ws.Range("D44").Value = "0.0₃0686"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.34"
$ws.Range("D50").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -4.15%  "
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("E31").Value = "  -3.67%  "
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("E35").Value = "  -15.14%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  -11.62%  "
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("E46").Value = "  -5.32%  "
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("E51").Value = "  +0.03%  "
